$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14 ("get avatar for each housemate ...") moves from IN PROGRESS -> DONE.
# Copy the formatting already used for "DONE" cells (e.g. C2) onto C14,
# then set its text to the already-existing "DONE" shared string.
$ws.Range("C2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 15 ("display room avatar ...") moves from NOT STARTED -> IN PROGRESS.
# Copy the formatting already used for the bordered "IN PROGRESS" look
# (same Neutral+border style as the Sprint 4 label cells) onto C15,
# then set its text to the already-existing "IN PROGRESS" shared string.
$ws.Range("A12").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C14").Value = "DONE"
$ws.Range("C15").Value = "IN PROGRESS"

# Update the saved cursor/selection position on the sheet.
$ws.Range("D13").Select()
